$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - F column "想去人数" (interested count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3098
$ws1.Range("F3").Value = 732
$ws1.Range("F4").Value = 112
$ws1.Range("F5").Value = 6826
$ws1.Range("F6").Value = 1872
$ws1.Range("F7").Value = 4
$ws1.Range("F8").Value = 57
$ws1.Range("F12").Value = 6
$ws1.Range("F14").Value = 158

# Sheet "全部类型" (All Types) - same underlying rows, shifted by one extra row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3098
$ws4.Range("F4").Value = 732
$ws4.Range("F5").Value = 112
$ws4.Range("F6").Value = 6826
$ws4.Range("F7").Value = 1872
$ws4.Range("F8").Value = 4
$ws4.Range("F9").Value = 57
$ws4.Range("F13").Value = 6
$ws4.Range("F15").Value = 158
